$d = $word.ActiveDocument

# 1) "in seguito il team ... View-Model-Controller." -> add comma after "in seguito" and
#    swap "View-Model-Controller" to "Model-View-Controller"
[void]$d.Content.Find.Execute(
    "in seguito il team ha ritenuto più opportuno optare per un altro stile architettonico: una variante dello stile View-Model-Controller.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "in seguito, il team ha ritenuto più opportuno optare per un altro stile architettonico: una variante dello stile Model-View-Controller.",
    2)

# 2) Remove the stray "_GoBack" bookmark left over from the previous save.
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}
